$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-6.19%"
$ws.Range("D3").Value = "'35.03"
$ws.Range("E3").Value = "'-3.31%"
$ws.Range("D4").Value = "'5.030"
$ws.Range("E4").Value = "'-1.87%"
$ws.Range("D5").Value = "'0.07907"
$ws.Range("E5").Value = "'-2.30%"
$ws.Range("D6").Value = "'1.942"
$ws.Range("E6").Value = "'-9.97%"
$ws.Range("D7").Value = "'7.748"
$ws.Range("E7").Value = "'-3.27%"
$ws.Range("D8").Value = "'4.020"
$ws.Range("E8").Value = "'-2.90%"
$ws.Range("E9").Value = "'5.97%"
$ws.Range("D10").Value = "'0.9236"
$ws.Range("E10").Value = "'-0.43%"
$ws.Range("D11").Value = "'0.1184"
$ws.Range("E11").Value = "'18.26%"
$ws.Range("D12").Value = "'0.1838"
$ws.Range("E12").Value = "'-2.91%"
$ws.Range("D13").Value = "'0.09318"
$ws.Range("E13").Value = "'1.33%"
$ws.Range("D14").Value = "'0.03530"
$ws.Range("E14").Value = "'-1.77%"
$ws.Range("D15").Value = "'0.09863"
$ws.Range("E15").Value = "'-0.77%"
$ws.Range("D16").Value = "'0.001387"
$ws.Range("E16").Value = "'-3.41%"
$ws.Range("D17").Value = "'0.005902"
$ws.Range("E17").Value = "'3.88%"
$ws.Range("D18").Value = "'3.505"
$ws.Range("E18").Value = "'1.37%"
$ws.Range("D19").Value = "'0.3443"
$ws.Range("E19").Value = "'2.10%"
$ws.Range("D20").Value = "'0.1308"
$ws.Range("E20").Value = "'-1.67%"
$ws.Range("D21").Value = "'5.046"
$ws.Range("E21").Value = "'-0.20%"
$ws.Range("D22").Value = "'0.2398"
$ws.Range("E22").Value = "'8.81%"
$ws.Range("D23").Value = "'0.04498"
$ws.Range("E23").Value = "'-2.23%"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'-2.21%"
$ws.Range("D25").Value = "'0.004568"
$ws.Range("E25").Value = "'-3.54%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'-3.85%"
$ws.Range("E27").Value = "'-6.89%"
$ws.Range("D39").Value = "'0.01901"
$ws.Range("E39").Value = "'-6.20%"
$ws.Range("D40").Value = "'0.04704"
$ws.Range("E40").Value = "'-5.69%"
$ws.Range("D41").Value = "'0.007594"
$ws.Range("E41").Value = "'-2.75%"
$ws.Range("D42").Value = "'0.009555"
$ws.Range("E42").Value = "'22.26%"
$ws.Range("D43").Value = "'0.1324"
$ws.Range("E43").Value = "'-5.45%"
$ws.Range("D44").Value = "'0.002111"
$ws.Range("E44").Value = "'1.41%"
$ws.Range("D45").Value = "'0.01118"
$ws.Range("E45").Value = "'-7.83%"
$ws.Range("D46").Value = "'0.00006008"
$ws.Range("E46").Value = "'-6.69%"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("E49").Value = "'-31.40%"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.00%"
